# Applies the "first pass of new method" edit:
# For each affected block of two header rows (row N / row N+1):
#   - Row N previously had its text only in column B; move it to column A
#     (clear column B).
#   - Row N+1 had column A holding the "~UC_..." marker and column B holding
#     the "~UC_Sets: R_*: AllRegions" text; swap them so column A gets the
#     AllRegions text and column B gets the marker text.
#
# This mirrors the exact cell moves shown in the OOXML diff for each of the
# affected worksheets/rows.

function Apply-HeaderSwap($ws, $row) {
    $nextRow = $row + 1

    # Current values before the edit.
    $bText = $ws.Cells.Item($row, 2).Value()         # e.g. "~UC_Sets: T_S: "
    $aNext = $ws.Cells.Item($nextRow, 1).Value()      # e.g. "~UC_T" / "~TFM_INS" / "~TFM_UPD"
    $bNext = $ws.Cells.Item($nextRow, 2).Value()      # e.g. "~UC_Sets: R_S: AllRegions"

    # Row N: move the text from B to A, clear B.
    $ws.Cells.Item($row, 1).Value = $bText
    $ws.Cells.Item($row, 2).ClearContents()

    # Row N+1: swap A and B contents.
    $ws.Cells.Item($nextRow, 1).Value = $bNext
    $ws.Cells.Item($nextRow, 2).Value = $aNext
}

$wb = $excel.ActiveWorkbook

# Sheet "Cars" - two header blocks (rows 1-2 and 7-8)
$ws1 = $wb.Worksheets.Item("Cars")
Apply-HeaderSwap $ws1 1
Apply-HeaderSwap $ws1 7

# Sheet "Cars_2020" - one header block (rows 1-2)
$ws2 = $wb.Worksheets.Item("Cars_2020")
Apply-HeaderSwap $ws2 1

# Sheet "CCS+h2" - one header block (rows 1-2)
$ws3 = $wb.Worksheets.Item("CCS+h2")
Apply-HeaderSwap $ws3 1

# Sheet "CH_RH" - one header block (rows 1-2)
$ws4 = $wb.Worksheets.Item("CH_RH")
Apply-HeaderSwap $ws4 1

# Sheet "IND_fuels" - one header block (rows 1-2)
$ws5 = $wb.Worksheets.Item("IND_fuels")
Apply-HeaderSwap $ws5 1

# Sheet "Power_sector" - three header blocks (rows 1-2, 10-11, 17-18)
$ws6 = $wb.Worksheets.Item("Power_sector")
Apply-HeaderSwap $ws6 1
Apply-HeaderSwap $ws6 10
Apply-HeaderSwap $ws6 17

# Sheet "Thermal_gencap" - two header blocks (rows 1-2 and 7-8)
$ws7 = $wb.Worksheets.Item("Thermal_gencap")
Apply-HeaderSwap $ws7 1
Apply-HeaderSwap $ws7 7

# Sheet "TRA_Policy" - one header block (rows 1-2)
$ws8 = $wb.Worksheets.Item("TRA_Policy")
Apply-HeaderSwap $ws8 1
